$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, pushing existing rows 42-47 down to 43-48.
$ws.Rows.Item(42).Insert()

# Populate the new row 42 with the new weekly price record.
$ws.Cells.Item(42, 1).Value2 = 11
$ws.Cells.Item(42, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(42, 3).Value2 = "Bíobío"
$ws.Cells.Item(42, 4).Value2 = 45034
$ws.Cells.Item(42, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(42, 5).Value2 = 8
$ws.Cells.Item(42, 6).Value2 = 100114007
$ws.Cells.Item(42, 7).Value2 = "Jengibre"
$ws.Cells.Item(42, 8).Value2 = "Sin especificar"
$ws.Cells.Item(42, 9).Value2 = "Primera"
$ws.Cells.Item(42, 10).Value2 = 50
$ws.Cells.Item(42, 11).Value2 = 15000
$ws.Cells.Item(42, 12).Value2 = 16000
$ws.Cells.Item(42, 13).Value2 = 15600
$ws.Cells.Item(42, 14).Value2 = "$/caja 13 kilos"
$ws.Cells.Item(42, 15).Value2 = "Perú"
$ws.Cells.Item(42, 16).Value2 = 1200
$ws.Cells.Item(42, 17).Value2 = 13
$ws.Cells.Item(42, 18).Value2 = "Hortaliza"
